# Generate Report for Handoff
# Re-stamps the "Latest Handoff Date/Datetime" for every row that is still
# pending a handoff (status "Handback transform failed" or "Ready for
# handoff") with the current run's single handoff timestamp, per sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Overview sheet: column D ("Latest Handoff Date") for rows 7, 10-16
$overviewRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $overviewRows) {
    $ws1.Cells.Item($r, 4).Value = "2016-20-12 12:20:58"
}

# zh-cn sheet: column E ("Latest Handoff Datetime") for rows 7, 10-16
$langRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $langRows) {
    $ws2.Cells.Item($r, 5).Value = "2016-03-12 12:20:54"
}

# de-de sheet: column E ("Latest Handoff Datetime") for rows 7, 10-16
foreach ($r in $langRows) {
    $ws3.Cells.Item($r, 5).Value = "2016-03-12 12:20:58"
}
